$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected columns (B-E) keep their existing text format so that
# numeric-looking strings (prices) are not coerced into actual numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.507.43"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "3.690.06"
$ws.Range("E3").Value = "  +4.48%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "419.38"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").Value = "129.93"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("D7").Value = "3.682.12"
$ws.Range("E7").Value = "  +4.51%  "
$ws.Range("D8").Value = "0.643"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").Value = "0.181"
$ws.Range("E11").Value = "  +8.69%  "
$ws.Range("D12").Value = "0.0000395"
$ws.Range("E12").Value = "  +45.81%  "
$ws.Range("D13").Value = "43.10"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  +5.49%  "
$ws.Range("D15").Value = "4.277.62"
$ws.Range("E15").Value = "  +4.64%  "
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "20.58"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "3.694.35"
$ws.Range("E18").Value = "  +4.32%  "
$ws.Range("D19").Value = "13.35"
$ws.Range("E19").Value = "  +6.07%  "
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "66.557.79"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").Value = "444.42"
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").Value = "16.47"
$ws.Range("E23").Value = "  +23.09%  "
$ws.Range("D24").Value = "90.00"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("D26").Value = "37.40"
$ws.Range("E26").Value = "  +8.51%  "
$ws.Range("D27").Value = "10.21"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").Value = "5.11"
$ws.Range("E29").Value = "  +5.93%  "
$ws.Range("E30").Value = "  +8.91%  "
$ws.Range("D31").Value = "12.72"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("D35").Value = "41.66"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("D36").Value = "57.35"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "0.0493"
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").Value = "3.09"
$ws.Range("E39").Value = "  +33.47%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0730"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").Value = "0.150"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("D42").Value = "29.38"
$ws.Range("E42").Value = "  +33.32%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").Value = "148.66"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("E46").Value = "  +4.49%  "
$ws.Range("D47").Value = "2.67"
$ws.Range("E47").Value = "  -4.53%  "
$ws.Range("D48").Value = "2.89"
$ws.Range("E48").Value = "  -7.69%  "
$ws.Range("D49").Value = "4.36"
$ws.Range("E49").Value = "  -4.72%  "
$ws.Range("D50").Value = "0.306"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("E51").Value = "  +12.02%  "
